$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell, copying the header style (bold + border + center/top align)
# from the existing "panel" header cell (E1) so it reuses the same style index.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

$timeTaken = @(
    "2021-10-05 13:42:28.059861",
    "2021-10-05 13:42:28.059874",
    "2021-10-05 13:42:28.059878",
    "2021-10-05 13:42:28.059881",
    "2021-10-05 13:42:28.059884",
    "2021-10-05 13:42:28.059887",
    "2021-10-05 13:42:28.059890",
    "2021-10-05 13:42:28.059893",
    "2021-10-05 13:42:28.059897",
    "2021-10-05 13:42:28.059900",
    "2021-10-05 13:42:28.059903",
    "2021-10-05 13:42:28.059905",
    "2021-10-05 13:42:28.059908",
    "2021-10-05 13:42:28.059911",
    "2021-10-05 13:42:28.059914",
    "2021-10-05 13:42:28.059917",
    "2021-10-05 13:42:28.059921",
    "2021-10-05 13:42:28.059924",
    "2021-10-05 13:42:28.059927",
    "2021-10-05 13:42:28.059930",
    "2021-10-05 13:42:28.059933",
    "2021-10-05 13:42:28.059936",
    "2021-10-05 13:42:28.059939",
    "2021-10-05 13:42:28.059942",
    "2021-10-05 13:42:28.059945"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
